$d = $word.ActiveDocument

# 1. Trim the leading question off the first run, keep "ezEML's "
$d.Content.Find.Execute(
    "How do you know when you" + [char]8217 + "ve filled in all of the required and recommended values? ezEML" + [char]8217 + "s ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ezEML" + [char]8217 + "s ",
    2
)

# 2. Expand the trailing run's text
$d.Content.Find.Execute(
    " feature will tell you.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " feature helps you determine when all of the required and recommended values have been completed.",
    2
)

# 3. Insert a new italic paragraph "…to be continued…" right after the blank
#    paragraph that follows the "Check" paragraph (and before the next two
#    blank paragraphs).
$checkPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*ezEML*Check feature*") {
        $checkPara = $p
        break
    }
}

$blankAfterCheck = $checkPara.Next()
$insertionRange = $blankAfterCheck.Range
$insertionRange.Collapse(0)
$newPara = $insertionRange.InsertParagraphAfter()

# Re-fetch the freshly minted (still empty) paragraph and stamp it with the
# exact OOXML we need (both the paragraph-mark rPr and run rPr get <w:i/>
# and <w:iCs/>) via InsertXML so formatting round-trips exactly.
$newBlank = $blankAfterCheck.Next()
$targetRange = $newBlank.Range
$targetRange.Collapse(1)

$ellipsis = [char]8230
$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>' + $ellipsis + 'to be continued' + $ellipsis + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$targetRange.InsertXML($xmlFrag)

Write-Host "Edit complete"
